# Apply the commit: insert a new "Coliflor" price record row right after
# row 423 (pushing the existing rows 424..474 down to 425..475), and give
# row 423 brand-new data (date 45142, vol 5000, prices 600).
#
# The new row 424 is first created as an exact duplicate of the old row
# 423 (via copy/paste of the whole row, including formatting), matching
# every row in this block (cols A,B,C,E,F,G,H,N,Q,R are constant across
# the whole table), then row 423 is overwritten in place with the new
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at 424 - this shifts old rows 424..474 down to 425..475.
$ws.Rows(424).Insert()

# 2) Duplicate row 423 (still holding the old data) into the newly inserted
#    row 424, so every column (incl. formatting/style) lines up correctly.
$ws.Range("A423:R423").Copy()
$ws.Range("A424").PasteSpecial()

# 3) Overwrite row 423 with its new values per the diff.
$ws.Range("D423").Value = 45142
$ws.Range("J423").Value = 5000
$ws.Range("K423").Value = 600
$ws.Range("L423").Value = 600
$ws.Range("M423").Value = 600
$ws.Range("O423").Value = "Región del Maule"
$ws.Range("P423").Value = 600
